$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B (old B..E become E..H)
$ws.Columns("B:D").Insert()

# New header row values for the inserted week columns + re-affirm the rest
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"
$ws.Range("B1").Value = "Jun_27"

# Fill the new columns (B,C,D) for existing analyst rows 2-27 with "UN"
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Add two new analyst rows for the new group
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
